# Update the high score name in row 10 from "Yo Dawg Crilla" to "Yeet"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Yeet"
